$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '91.045.75'
$ws.Range("E2").Value = '  +4.37%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '3.190.49'
$ws.Range("E3").Value = '  +1.36%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '218.87'
$ws.Range("E5").Value = '  +6.10%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '634.02'
$ws.Range("E6").Value = '  +5.13%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.401'
$ws.Range("E7").Value = '  +4.78%  '

$ws.Range("E8").Value = '  +7.45%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '3.188.39'
$ws.Range("E10").Value = '  +1.38%  '

$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.572'
$ws.Range("E11").Value = '  +7.86%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '0.181'
$ws.Range("E12").Value = '  +3.53%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '0.0000258'
$ws.Range("E13").Value = '  +6.81%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '5.40'
$ws.Range("E14").Value = '  +3.40%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '90.759.59'
$ws.Range("E15").Value = '  +4.22%  '

$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '3.782.77'
$ws.Range("E16").Value = '  +1.27%  '

$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '33.08'
$ws.Range("E17").Value = '  +3.82%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '3.203.60'
$ws.Range("E18").Value = '  +1.31%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '0.0000226'
$ws.Range("E19").Value = '  +75.24%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '3.32'
$ws.Range("E20").Value = '  +5.21%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '439.39'
$ws.Range("E21").Value = '  +6.88%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '13.39'
$ws.Range("E22").Value = '  +0.74%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '8.57'
$ws.Range("E23").Value = '  +1.88%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '5.03'
$ws.Range("E24").Value = '  +0.68%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '5.29'
$ws.Range("E25").Value = '  +3.68%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '11.77'
$ws.Range("E26").Value = '  -0.62%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '80.83'
$ws.Range("E27").Value = '  +10.97%  '

$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '3.368.75'
$ws.Range("E28").Value = '  +1.25%  '

$ws.Range("E29").Value = '  +0.10%  '

$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.69%  '

$ws.Range("B31").Value = 'Cronos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '0.158'
$ws.Range("E31").Value = '  -1.60%  '

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '4.12'
$ws.Range("E32").Value = '  +37.82%  '

$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '8.39'
$ws.Range("E33").Value = '  +3.15%  '

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '529.75'
$ws.Range("E34").Value = '  -2.00%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '7.01'
$ws.Range("E35").Value = '  +5.94%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '1.90'
$ws.Range("E36").Value = '  +3.63%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '1.29'
$ws.Range("E37").Value = '  -0.29%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '22.45'
$ws.Range("E38").Value = '  +3.63%  '

$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.08%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '0.127'
$ws.Range("E41").Value = '  -2.27%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '1.94'
$ws.Range("E42").Value = '  +2.69%  '

$ws.Range("E43").Value = '  +0.01%  '

$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '0.373'
$ws.Range("E44").Value = '  +2.13%  '

$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '147.34'
$ws.Range("E45").Value = '  -1.58%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '44.23'
$ws.Range("E46").Value = '  +2.67%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '172.50'
$ws.Range("E47").Value = '  +0.58%  '

$ws.Range("E48").Value = '  +2.22%  '

$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '0.747'
$ws.Range("E49").Value = '  +8.71%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '24.90'
$ws.Range("E50").Value = '  +6.68%  '

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '1.23'
$ws.Range("E51").Value = '  +1.38%  '
